$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27:C27").NumberFormat = "@"

$ws.Range("A27").Value = "2025-09-19"
$ws.Range("B27").Value = "15:20:55"
$ws.Range("C27").Value = "1.00 EUR = 1,749.0925"
